# "Historial de Comits" - add the second commit-log entry ("El error del
# vacio persiste") to the tracking sheet, matching the author's original
# table of columns: Fecha / Integrante / Titulo del Comit / Descripcion del
# Comit, with centered/wrapped formatting and a hyperlink on the first
# commit's title.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlCenter = -4108
$xlPasteFormats = -4122
$xlPortrait = 1

# --- Column widths (B/C/D) ---------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.833333333333334
$ws.Columns.Item(3).ColumnWidth = 26.333333333333332
$ws.Columns.Item(4).ColumnWidth = 51.666666666666664

# --- Base style: centered horizontally + vertically (rows 1-21) --------
$ws.Range("A1").VerticalAlignment = $xlCenter
$ws.Range("A1").HorizontalAlignment = $xlCenter
$ws.Range("A1").Copy()
$ws.Range("A1:D21").PasteSpecial($xlPasteFormats)

# --- Header row (row 2) --------------------------------------------------
$ws.Range("A2").Value = "Fecha "
$ws.Range("B2").Value = "Integrante "
$ws.Range("C2").Value = "Titulo del Comit"
$ws.Range("D2").Value = "Descripcion del Comit"

# --- Date style (fecha column, rows 3 & 4) -------------------------------
$ws.Range("A3").NumberFormat = "d-mmm"
$ws.Range("A3").VerticalAlignment = $xlCenter
$ws.Range("A3").HorizontalAlignment = $xlCenter
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)

# --- Wrapped + centered style (Integrante/Descripcion cells) ------------
$ws.Range("D3").WrapText = $true
$ws.Range("D3").VerticalAlignment = $xlCenter
$ws.Range("D3").HorizontalAlignment = $xlCenter
$ws.Range("D3").Copy()
$ws.Range("B3").PasteSpecial($xlPasteFormats)
$ws.Range("B4").PasteSpecial($xlPasteFormats)
$ws.Range("D4").PasteSpecial($xlPasteFormats)

# --- First commit title: hyperlink (sets display text) + real cell text +
#     wrapped/centered style. Hyperlinks.Add's TextToDisplay overwrites the
#     cell, so set the real title afterwards to restore it.
$ws.Hyperlinks.Add($ws.Cells.Item(3, 3), "https://github.com/ipopotamo/I.R.O.J.I/commit/20f1e72b5ff7299972145393c3bbdd169b4b2609", "", "", "https://github.com/ipopotamo/I.R.O.J.I/commit/20f1e72b5ff7299972145393c3bbdd169b4b2609")
$ws.Range("C3").Value = "Brumer arregla bug al vacio Apollyon Triste"
$ws.Range("C3").WrapText = $true
$ws.Range("C3").VerticalAlignment = $xlCenter
$ws.Range("C3").HorizontalAlignment = $xlCenter

# --- Remaining cell text (order matters for shared-string indices) ------
$ws.Range("D3").Value = "Se arregla el error en el que 2 salas chocan y generan un vacio"
$ws.Range("B3").Value = "Octavio Lucardi Fierro"
$ws.Range("B4").Value = "Octavio Lucardi Fierro"
$ws.Range("C4").Value = "El error del vacio persiste"
$ws.Range("D4").Value = "Al actualizar el repositorio y correr el codigo, este ya no reconoce la librería de Unity y deja de funcionar todo"
$ws.Range("A3").Value = 44721
$ws.Range("A4").Value = 44721

# --- Row heights ----------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 50.25
$ws.Rows.Item(4).RowHeight = 45

# --- Page setup + selection ------------------------------------------------
$ws.PageSetup.Orientation = $xlPortrait
$ws.Range("C4").Select()
